$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

Set-TextValue ($ws.Range("D2")) '289.93'
Set-TextValue ($ws.Range("E2")) '-3.86%'
Set-TextValue ($ws.Range("D3")) '30.58'
Set-TextValue ($ws.Range("E3")) '-6.28%'
Set-TextValue ($ws.Range("D4")) '4.917'
Set-TextValue ($ws.Range("E4")) '-3.20%'
Set-TextValue ($ws.Range("D5")) '0.07251'
Set-TextValue ($ws.Range("E5")) '-6.10%'
Set-TextValue ($ws.Range("D6")) '1.803'
Set-TextValue ($ws.Range("E6")) '-11.53%'
Set-TextValue ($ws.Range("D7")) '7.635'
Set-TextValue ($ws.Range("E7")) '-3.53%'
Set-TextValue ($ws.Range("D8")) '3.698'
Set-TextValue ($ws.Range("E8")) '-2.74%'
Set-TextValue ($ws.Range("D9")) '0.9011'
Set-TextValue ($ws.Range("E9")) '-2.73%'
Set-TextValue ($ws.Range("D10")) '0.1679'
Set-TextValue ($ws.Range("E10")) '-4.76%'
Set-TextValue ($ws.Range("D11")) '0.08028'
Set-TextValue ($ws.Range("E11")) '-1.60%'
Set-TextValue ($ws.Range("D12")) '0.08100'
Set-TextValue ($ws.Range("E12")) '-5.78%'
Set-TextValue ($ws.Range("D13")) '0.03053'
Set-TextValue ($ws.Range("E13")) '0.06%'
Set-TextValue ($ws.Range("E14")) '0.29%'
Set-TextValue ($ws.Range("D15")) '0.001496'
Set-TextValue ($ws.Range("E15")) '-1.81%'
Set-TextValue ($ws.Range("D16")) '0.005699'
Set-TextValue ($ws.Range("E16")) '-3.39%'
Set-TextValue ($ws.Range("D17")) '3.480'
Set-TextValue ($ws.Range("E17")) '0.12%'
Set-TextValue ($ws.Range("D18")) '2.075'
Set-TextValue ($ws.Range("E18")) '-3.80%'
Set-TextValue ($ws.Range("D19")) '0.3319'
Set-TextValue ($ws.Range("E19")) '-0.44%'
Set-TextValue ($ws.Range("D20")) '0.1302'
Set-TextValue ($ws.Range("E20")) '-1.91%'
Set-TextValue ($ws.Range("D21")) '3.967'
Set-TextValue ($ws.Range("E21")) '-10.01%'
Set-TextValue ($ws.Range("E22")) '9.72%'
Set-TextValue ($ws.Range("D23")) '0.04506'
Set-TextValue ($ws.Range("E23")) '-0.85%'
Set-TextValue ($ws.Range("E24")) '-1.48%'
Set-TextValue ($ws.Range("D25")) '0.004437'
Set-TextValue ($ws.Range("E25")) '6.82%'
Set-TextValue ($ws.Range("D26")) '0.0001300'
Set-TextValue ($ws.Range("E26")) '3.77%'
Set-TextValue ($ws.Range("D27")) '0.0003387'
Set-TextValue ($ws.Range("E27")) '-95.48%'
Set-TextValue ($ws.Range("D39")) '0.01584'
Set-TextValue ($ws.Range("E39")) '-8.70%'
Set-TextValue ($ws.Range("D40")) '0.04353'
Set-TextValue ($ws.Range("E40")) '-7.46%'
Set-TextValue ($ws.Range("D41")) '0.007349'
Set-TextValue ($ws.Range("E41")) '-1.88%'
Set-TextValue ($ws.Range("D42")) '0.01002'
Set-TextValue ($ws.Range("D43")) '0.1315'
Set-TextValue ($ws.Range("E43")) '-3.63%'
Set-TextValue ($ws.Range("D44")) '0.002004'
Set-TextValue ($ws.Range("E44")) '-14.17%'
Set-TextValue ($ws.Range("D45")) '0.009454'
Set-TextValue ($ws.Range("E45")) '-9.13%'
Set-TextValue ($ws.Range("D46")) '0.00005856'
Set-TextValue ($ws.Range("E46")) '-5.05%'
Set-TextValue ($ws.Range("D47")) '0.00000000748'
Set-TextValue ($ws.Range("E47")) '-0.29%'
Set-TextValue ($ws.Range("D48")) '2.251'
Set-TextValue ($ws.Range("E48")) '25.03%'
Set-TextValue ($ws.Range("D49")) '0.002893'
Set-TextValue ($ws.Range("E49")) '-3.41%'
Set-TextValue ($ws.Range("D50")) '0.00002096'
Set-TextValue ($ws.Range("E50")) '-0.29%'
Set-TextValue ($ws.Range("D51")) '0.0001996'
Set-TextValue ($ws.Range("E51")) '-0.29%'
